$wb = $excel.ActiveWorkbook

$newName = "dfc3b845-555e-4fc7-801d-c410bc78d5f4.md"
$epoch   = "0001-01-01 00:00:00"
$status  = "Handoff failed"
$reason  = "Ignored"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh       = $wb.Worksheets.Item("zh-cn")
$wsDe       = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Re-point the "handoff file" hyperlink (column A, row 2) on every sheet to
# the new handback report file, and refresh the status text that goes with
# a failed handoff.
# ---------------------------------------------------------------------------

$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/028c8c5120be7fe88fba16b661b7860623a663c4/e2e/$newName",
    [Type]::Missing,
    [Type]::Missing,
    $newName)
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/028c8c5120be7fe88fba16b661b7860623a663c4/.localization-config",
    [Type]::Missing,
    [Type]::Missing,
    ".localization-config")

$wsOverview.Range("B2").Value = $status
$wsOverview.Range("C2").Value = $status

$wsZh.Range("A1").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add(
    $wsZh.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/028c8c5120be7fe88fba16b661b7860623a663c4/e2e/$newName",
    [Type]::Missing,
    [Type]::Missing,
    $newName)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/028c8c5120be7fe88fba16b661b7860623a663c4/.localization-config",
    [Type]::Missing,
    [Type]::Missing,
    ".localization-config")

$wsDe.Range("A1").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add(
    $wsDe.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/028c8c5120be7fe88fba16b661b7860623a663c4/e2e/$newName",
    [Type]::Missing,
    [Type]::Missing,
    $newName)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/028c8c5120be7fe88fba16b661b7860623a663c4/.localization-config",
    [Type]::Missing,
    [Type]::Missing,
    ".localization-config")

# ---------------------------------------------------------------------------
# zh-cn detail sheet: the handoff attempt is now failed/ignored, so the
# "Latest Handoff File" link is removed and the datetimes reset to the zero
# value, matching a dependency that was never actually handed off.
# ---------------------------------------------------------------------------
$wsZh.Range("B2").Value = $status
$wsZh.Range("C2").Clear()
$wsZh.Range("D2").Value = $epoch
$wsZh.Range("G2").Value = $epoch
$wsZh.Range("H2").Value = $reason

# ---------------------------------------------------------------------------
# de-de detail sheet: identical treatment.
# ---------------------------------------------------------------------------
$wsDe.Range("B2").Value = $status
$wsDe.Range("C2").Clear()
$wsDe.Range("D2").Value = $epoch
$wsDe.Range("G2").Value = $epoch
$wsDe.Range("H2").Value = $reason

Write-Host "Report regenerated for handoff failure."
